# "unify the conception of DataNode, DataTable, Entity."
#
# Rename the two sheets to their new, unified names:
#   Property1      -> DataNode
#   Record_Station -> DataTable
# then make DataTable the active/selected tab (it previously was
# Property1/DataNode) with its own selected cell, and nudge a couple of
# DataNode row heights that changed with the new header wrapping.

$wb = $excel.ActiveWorkbook

$wsDataNode  = $wb.Worksheets.Item("Property1")
$wsDataTable = $wb.Worksheets.Item("Record_Station")

$wsDataNode.Name  = "DataNode"
$wsDataTable.Name = "DataTable"

# Header row now wraps onto two lines; the long-description row shrank.
$wsDataNode.Rows.Item(1).RowHeight = 27
$wsDataNode.Rows.Item(8).RowHeight = 54

# DataTable becomes the visible/active sheet, selection moved to H32.
$wsDataTable.Activate()
$wsDataTable.Range("H32").Select()
